# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and apply the OKB / dogwifhat row reordering (rows 44-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.770.90'
$ws.Range('E2').Value = '  -3.53%  '
$ws.Range('D3').Value = '2.483.94'
$ws.Range('E3').Value = '  -6.23%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '556.64'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -4.36%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '148.16'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -5.38%  '
$ws.Range('E7').Value = '  -0.07%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.601'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -3.58%  '
$ws.Range('D9').Value = '2.480.11'
$ws.Range('E9').Value = '  -6.29%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.109'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -8.40%  '
$ws.Range('E11').Value = '  -5.47%  '
$ws.Range('E12').Value = '  -1.46%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.360'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -6.65%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '26.56'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -7.14%  '
$ws.Range('D15').Value = '2.927.84'
$ws.Range('E15').Value = '  -6.32%  '
$ws.Range('E16').Value = '  -8.48%  '
$ws.Range('D17').Value = '61.576.39'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').Value = '2.481.57'
$ws.Range('E18').Value = '  -6.30%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.27'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -7.72%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.19'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -7.34%  '
$ws.Range('E21').Value = '  -6.71%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '323.17'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -6.80%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  +1.91%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '64.50'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -5.42%  '
$ws.Range('E26').Value = '  -9.08%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '568.83'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('D28').Value = '2.605.56'
$ws.Range('E28').Value = '  -6.19%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.52'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -6.29%  '
$ws.Range('E30').Value = '  -0.18%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.85'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -4.73%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '8.36'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -10.23%  '
$ws.Range('E33').Value = '  -6.64%  '
$ws.Range('E34').Value = '  -6.09%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.60'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -8.30%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.99'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -9.93%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.95'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -10.57%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('E39').Value = '  -4.74%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '18.66'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -5.78%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '146.44'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -3.16%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.78'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -8.00%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.46'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -4.48%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '40.62'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -3.14%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '149.22'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -8.72%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.67'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -6.53%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '22.12'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -9.61%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0544'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -7.98%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.599'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -5.96%  '
$ws.Range('E51').Value = '  -5.52%  '
